# "added JP lastName on file"
#
# The paragraph that lists the authors currently ends with a red,
# spell-check-flagged placeholder "(Jp ponga su apellido xD)" ("Jp put your
# last name lol"). Replace that placeholder with the real last name
# ("Betancourt Maldonado") in plain (non-highlighted) formatting, and move
# the document's "_GoBack" last-edit bookmark to sit right after the newly
# typed name (its previous location, after "permitir " further down, is
# vacated).

$d = $word.ActiveDocument

# Locate the paragraph containing the author list / placeholder text.
$authorsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Bonilla*(Jp*") {
        $authorsPara = $p
        break
    }
}

$paraStart = $authorsPara.Range.Start
$paraText = $authorsPara.Range.Text

# Position right before the opening "(" of the "(Jp ponga su apellido xD)"
# placeholder -- i.e. right after "..., Juan Pablo ".
$placeholderOffset = $paraText.IndexOf("(Jp ponga su apellido xD)")
$insertPos = $paraStart + $placeholderOffset

# Type the real last name in directly (inherits the plain, uncoloured
# formatting of the text immediately before it, unlike the red placeholder).
$insRange = $d.Range($insertPos, $insertPos)
$insRange.InsertAfter("Betancourt Maldonado")

# The new text now sits just before the still-present placeholder, e.g.
# "...Juan Pablo Betancourt Maldonado(Jp ponga su apellido xD)". Re-home the
# "_GoBack" bookmark (collapsed) right after the name we just typed, while
# the placeholder text after it still keeps that position from being the
# paragraph's last character (adding a collapsed bookmark exactly at a
# paragraph's final character position is unreliable, so we do this before
# deleting the placeholder).
$bookmarkPos = $insertPos + ("Betancourt Maldonado").Length
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Now remove the placeholder text; the bookmark we just set stays put.
$placeholderLen = ("(Jp ponga su apellido xD)").Length
$d.Range($bookmarkPos, $bookmarkPos + $placeholderLen).Text = ""
